$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").ClearFormats()
$ws.Rows(2).Delete()
$ws.Range("A1").Value = 'questions = [
    {
        "title": "You have a raised button in the user interface (UI) of your app, as shown below. When the user taps on the button, you want to print \u201cContact Us\u201d in the debug console.Which code should you use?",
        "ques_type": 2,
        "options": [
            "RaisedButton(\n child: Text(''Contact us''),\n onPressed: () =&gt print(''Contact us'')\n),\n",
            "RaisedButton(\n child: Text(''Contact us''),\n onPressed: print(''Contact us'')\n),\n",
            "RaisedButton(\n child: Text(''Contact us''),\n onPressed: {print(''Contact us'')}\n),\n",
            "RaisedButton(\n child: Text(''Contact us''),\n onPressed: () =&gt return ''Contact us''\n),\n"
        ],
        "score": "RaisedButton(\n child: Text(''Contact us''),\n onPressed: () =&gt print(''Contact us'')\n),"
    },
    {
        "title": "In your app, you have two screens. You want to enable transitioning from the first screen to the second screen and back to the first screen without directly implementing any navigation code. Instead, you want to use a method of the Navigator object.Which method should you use to achieve this?",
        "ques_type": 2,
        "options": [
            "push()",
            "pop()",
            "pushReplacement()",
            "popAndPushNamed()"
        ],
        "score": "push()"
    },
    {
        "title": "You are a Flutter developer working on a mobile app that fetches data from a web service using the HTTP Dart package. You''ve just written the code snippet below, and your next step is to process the data.What action should you take?var myResult = await http.get(url)",
        "ques_type": 2,
        "options": [
            "Access the data from myResult using .body as it is an http.Response object.",
            "Convert myResult to a string using .toString() and process it.\u00a0",
            "Treat myResult as a Future object and use .then() to process it.",
            "Parse myResult as a JSON object using json.decode()."
        ],
        "score": "Access the data from myResult using .body as it is an http.Response object."
    },
    {
        "title": "You are a senior software engineer at a tech company. You''re mentoring a junior developer on a mobile app feature. The feature requires retrieving user profiles from an external web service asynchronously. You''re discussing function declarations for this task, particularly to get the name of the user.Which function declaration should you recommend?",
        "ques_type": 2,
        "options": [
            "String getName() async {}",
            "void getName() async {} ",
            "Future getName() async {} ",
            "Future&ltString&gt getName() async {}"
        ],
        "score": "String getName() async {}"
    }
]'
$ws.Rows(1).AutoFit()
